$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "66.158.50"
$cell.Style = "Normal"

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.313.81"
$cell.Style = "Normal"

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -0.29%  "
$cell.Style = "Normal"

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  +0.00%  "
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "585.32"
$cell.Style = "Normal"

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  +1.93%  "
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "180.84"
$cell.Style = "Normal"

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  +0.18%  "
$cell.Style = "Normal"

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  +5.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "3.317.28"
$cell.Style = "Normal"

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -0.14%  "
$cell.Style = "Normal"

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -1.32%  "
$cell.Style = "Normal"

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  +2.16%  "
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.400"
$cell.Style = "Normal"

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -0.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "3.899.24"
$cell.Style = "Normal"

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -0.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -2.55%  "
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "66.223.24"
$cell.Style = "Normal"

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -0.20%  "
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "26.56"
$cell.Style = "Normal"

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -0.39%  "
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "3.320.67"
$cell.Style = "Normal"

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  +0.24%  "
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.0000163"
$cell.Style = "Normal"

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -1.67%  "
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "423.27"
$cell.Style = "Normal"

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -3.07%  "
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "5.50"
$cell.Style = "Normal"

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -3.01%  "
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.09"
$cell.Style = "Normal"

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -3.21%  "
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "7.36"
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -2.85%  "
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "71.74"
$cell.Style = "Normal"

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -2.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -0.07%  "
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "5.67"
$cell.Style = "Normal"

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "3.474.21"
$cell.Style = "Normal"

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  +0.09%  "
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.512"
$cell.Style = "Normal"

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -1.18%  "
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.205"
$cell.Style = "Normal"

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.0000114"
$cell.Style = "Normal"

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -1.25%  "
$cell.Style = "Normal"

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "9.01"
$cell.Style = "Normal"

$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -0.34%  "
$cell.Style = "Normal"

$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.Style = "Normal"

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.92"
$cell.Style = "Normal"

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -1.49%  "
$cell.Style = "Normal"

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "22.34"
$cell.Style = "Normal"

$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -1.82%  "
$cell.Style = "Normal"

$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.Style = "Normal"

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.17"
$cell.Style = "Normal"

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -1.37%  "
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.61"
$cell.Style = "Normal"

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -1.84%  "
$cell.Style = "Normal"

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.18"
$cell.Style = "Normal"

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -2.62%  "
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "159.96"
$cell.Style = "Normal"

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -0.17%  "
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.43"
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -2.63%  "
$cell.Style = "Normal"

$cell = $ws.Range("B40")
$cell.NumberFormat = "@"
$cell.Value = "Maker"
$cell.Style = "Normal"

$cell = $ws.Range("C40")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.862.91"
$cell.Style = "Normal"

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +1.07%  "
$cell.Style = "Normal"

$cell = $ws.Range("B41")
$cell.NumberFormat = "@"
$cell.Value = "Stacks"
$cell.Style = "Normal"

$cell = $ws.Range("C41")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.80"
$cell.Style = "Normal"

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  +0.40%  "
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "26.36"
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -4.91%  "
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "4.32"
$cell.Style = "Normal"

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -2.23%  "
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.757"
$cell.Style = "Normal"

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -4.37%  "
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "39.73"
$cell.Style = "Normal"

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -1.24%  "
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0659"
$cell.Style = "Normal"

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -0.70%  "
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "5.95"
$cell.Style = "Normal"

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -3.68%  "
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.31"
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -1.39%  "
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "23.00"
$cell.Style = "Normal"

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -4.45%  "
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "309.66"
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -4.14%  "
$cell.Style = "Normal"

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  +0.42%  "
$cell.Style = "Normal"
